$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Cells.Item(2, 2).Value2 = 1.02
$ws.Cells.Item(2, 3).Value2 = 1.062627701646876
$ws.Cells.Item(2, 4).Value2 = 1.061933525267731
$ws.Cells.Item(2, 5).Value2 = 1.067255964106053
$ws.Cells.Item(2, 6).Value2 = 1.075901316290951
$ws.Cells.Item(2, 9).Value2 = 1.041617312011347
$ws.Cells.Item(2, 10).Value2 = 1.067597248363691
$ws.Cells.Item(2, 11).Value2 = 1.06465624103064
$ws.Cells.Item(2, 12).Value2 = 1.069964315849649
$ws.Cells.Item(2, 13).Value2 = 1.07858666363622
$ws.Cells.Item(2, 14).Value2 = 1.069113358911484

# Row 3
$ws.Cells.Item(3, 2).Value2 = 1.02
$ws.Cells.Item(3, 3).Value2 = 1.064136275271315
$ws.Cells.Item(3, 4).Value2 = 1.063072374536812
$ws.Cells.Item(3, 5).Value2 = 1.068576752876452
$ws.Cells.Item(3, 6).Value2 = 1.077254018091754
$ws.Cells.Item(3, 9).Value2 = 1.041919283356742
$ws.Cells.Item(3, 10).Value2 = 1.068757670267288
$ws.Cells.Item(3, 11).Value2 = 1.065608992134299
$ws.Cells.Item(3, 12).Value2 = 1.071099597479117
$ws.Cells.Item(3, 13).Value2 = 1.079755458041245
$ws.Cells.Item(3, 14).Value2 = 1.070275428747286

# Row 4
$ws.Cells.Item(4, 2).Value2 = 1.02
$ws.Cells.Item(4, 3).Value2 = 1.065111420953603
$ws.Cells.Item(4, 4).Value2 = 1.063808188716142
$ws.Cells.Item(4, 5).Value2 = 1.069430712599303
$ws.Cells.Item(4, 6).Value2 = 1.078128636270939
$ws.Cells.Item(4, 9).Value2 = 1.042112836090472
$ws.Cells.Item(4, 10).Value2 = 1.069507113538781
$ws.Cells.Item(4, 11).Value2 = 1.066223805368483
$ws.Cells.Item(4, 12).Value2 = 1.07183297391463
$ws.Cells.Item(4, 13).Value2 = 1.080510531494961
$ws.Cells.Item(4, 14).Value2 = 1.071025936314186

# Row 5
$ws.Cells.Item(5, 2).Value2 = 1.02
$ws.Cells.Item(5, 3).Value2 = 1.065521138650251
$ws.Cells.Item(5, 4).Value2 = 1.064117266121378
$ws.Cells.Item(5, 5).Value2 = 1.069789559413892
$ws.Cells.Item(5, 6).Value2 = 1.078496169956237
$ws.Cells.Item(5, 9).Value2 = 1.042193765800931
$ws.Cells.Item(5, 10).Value2 = 1.069821842388204
$ws.Cells.Item(5, 11).Value2 = 1.066481874192876
$ws.Cells.Item(5, 12).Value2 = 1.072140996001857
$ws.Cells.Item(5, 13).Value2 = 1.080827677869135
$ws.Cells.Item(5, 14).Value2 = 1.071341112114675

# Row 6
$ws.Cells.Item(6, 2).Value2 = 1.02
$ws.Cells.Item(6, 3).Value2 = 1.065589918544795
$ws.Cells.Item(6, 4).Value2 = 1.064169146508781
$ws.Cells.Item(6, 5).Value2 = 1.069849802252202
$ws.Cells.Item(6, 6).Value2 = 1.078557871478869
$ws.Cells.Item(6, 9).Value2 = 1.042207328500662
$ws.Cells.Item(6, 10).Value2 = 1.0698746671122
$ws.Cells.Item(6, 11).Value2 = 1.066525181811914
$ws.Cells.Item(6, 12).Value2 = 1.072192697421834
$ws.Cells.Item(6, 13).Value2 = 1.080880911464336
$ws.Cells.Item(6, 14).Value2 = 1.071394011855839

# Row 7
$ws.Cells.Item(7, 2).Value2 = 1.02
$ws.Cells.Item(7, 3).Value2 = 1.065116896531258
$ws.Cells.Item(7, 4).Value2 = 1.063812319635998
$ws.Cells.Item(7, 5).Value2 = 1.069435508142456
$ws.Cells.Item(7, 6).Value2 = 1.078133547880408
$ws.Cells.Item(7, 9).Value2 = 1.042113919203919
$ws.Cells.Item(7, 10).Value2 = 1.069511320279969
$ws.Cells.Item(7, 11).Value2 = 1.06622725525871
$ws.Cells.Item(7, 12).Value2 = 1.071837090852523
$ws.Cells.Item(7, 13).Value2 = 1.080514770340801
$ws.Cells.Item(7, 14).Value2 = 1.07103014902943

# Row 8
$ws.Cells.Item(8, 2).Value2 = 1.02
$ws.Cells.Item(8, 3).Value2 = 1.063137742571736
$ws.Cells.Item(8, 4).Value2 = 1.062318634272098
$ws.Cells.Item(8, 5).Value2 = 1.067702474258146
$ws.Cells.Item(8, 6).Value2 = 1.076358609821496
$ws.Cells.Item(8, 9).Value2 = 1.041719747038933
$ws.Cells.Item(8, 10).Value2 = 1.067989716682653
$ws.Cells.Item(8, 11).Value2 = 1.06497857787031
$ws.Cells.Item(8, 12).Value2 = 1.070348246184956
$ws.Cells.Item(8, 13).Value2 = 1.078981917274944
$ws.Cells.Item(8, 14).Value2 = 1.069506384580475

# Row 9
$ws.Cells.Item(9, 2).Value2 = 1.02
$ws.Cells.Item(9, 3).Value2 = 1.059642255979609
$ws.Cells.Item(9, 4).Value2 = 1.059677991081016
$ws.Cells.Item(9, 5).Value2 = 1.064643236755436
$ws.Cells.Item(9, 6).Value2 = 1.073225597630127
$ws.Cells.Item(9, 9).Value2 = 1.041010990326281
$ws.Cells.Item(9, 10).Value2 = 1.065297318782571
$ws.Cells.Item(9, 11).Value2 = 1.062765214767107
$ws.Cells.Item(9, 12).Value2 = 1.06771512303275
$ws.Cells.Item(9, 13).Value2 = 1.076271338325263
$ws.Cells.Item(9, 14).Value2 = 1.066810163166552

# Row 10
$ws.Cells.Item(10, 2).Value2 = 1.02
$ws.Cells.Item(10, 3).Value2 = 1.057306172573146
$ws.Cells.Item(10, 4).Value2 = 1.05791154930505
$ws.Cells.Item(10, 5).Value2 = 1.062599810819179
$ws.Cells.Item(10, 6).Value2 = 1.071133032723813
$ws.Cells.Item(10, 9).Value2 = 1.040528869642467
$ws.Cells.Item(10, 10).Value2 = 1.063494607441621
$ws.Cells.Item(10, 11).Value2 = 1.061280644816538
$ws.Cells.Item(10, 12).Value2 = 1.065952987975276
$ws.Cells.Item(10, 13).Value2 = 1.074457629074257
$ws.Cells.Item(10, 14).Value2 = 1.065004891768724

# Row 11
$ws.Cells.Item(11, 2).Value2 = 1.02
$ws.Cells.Item(11, 3).Value2 = 1.056293166249616
$ws.Cells.Item(11, 4).Value2 = 1.057145179032417
$ws.Cells.Item(11, 5).Value2 = 1.061713985001914
$ws.Cells.Item(11, 6).Value2 = 1.070225939602622
$ws.Cells.Item(11, 9).Value2 = 1.040317805754631
$ws.Cells.Item(11, 10).Value2 = 1.062712103009806
$ws.Cells.Item(11, 11).Value2 = 1.060635624151006
$ws.Cells.Item(11, 12).Value2 = 1.065188308915976
$ws.Cells.Item(11, 13).Value2 = 1.073670632879951
$ws.Cells.Item(11, 14).Value2 = 1.064221276090862

# Row 12
$ws.Cells.Item(12, 2).Value2 = 1.02
$ws.Cells.Item(12, 3).Value2 = 1.055916662358127
$ws.Cells.Item(12, 4).Value2 = 1.056860286379335
$ws.Cells.Item(12, 5).Value2 = 1.061384792447529
$ws.Cells.Item(12, 6).Value2 = 1.069888848947024
$ws.Cells.Item(12, 9).Value2 = 1.040239059571763
$ws.Cells.Item(12, 10).Value2 = 1.062421152068932
$ws.Cells.Item(12, 11).Value2 = 1.060395700841736
$ws.Cells.Item(12, 12).Value2 = 1.064904017858564
$ws.Cells.Item(12, 13).Value2 = 1.073378054478548
$ws.Cells.Item(12, 14).Value2 = 1.063929911966279

# Row 13
$ws.Cells.Item(13, 2).Value2 = 1.02
$ws.Cells.Item(13, 3).Value2 = 1.055997434176104
$ws.Cells.Item(13, 4).Value2 = 1.056921407244016
$ws.Cells.Item(13, 5).Value2 = 1.06145541258481
$ws.Cells.Item(13, 6).Value2 = 1.069961163185086
$ws.Cells.Item(13, 9).Value2 = 1.040255966655607
$ws.Cells.Item(13, 10).Value2 = 1.06248357542765
$ws.Cells.Item(13, 11).Value2 = 1.060447180390228
$ws.Cells.Item(13, 12).Value2 = 1.064965010890685
$ws.Cells.Item(13, 13).Value2 = 1.073440825072929
$ws.Cells.Item(13, 14).Value2 = 1.063992423973326

# Row 14
$ws.Cells.Item(14, 2).Value2 = 1.02
$ws.Cells.Item(14, 3).Value2 = 1.056262049022196
$ws.Cells.Item(14, 4).Value2 = 1.05712163440866
$ws.Cells.Item(14, 5).Value2 = 1.061686777087535
$ws.Cells.Item(14, 6).Value2 = 1.070198078797731
$ws.Cells.Item(14, 9).Value2 = 1.040311303672807
$ws.Cells.Item(14, 10).Value2 = 1.062688058944049
$ws.Cells.Item(14, 11).Value2 = 1.060615798853778
$ws.Cells.Item(14, 12).Value2 = 1.065164814565821
$ws.Cells.Item(14, 13).Value2 = 1.073646453449761
$ws.Cells.Item(14, 14).Value2 = 1.064197197879774

# Row 15
$ws.Cells.Item(15, 2).Value2 = 1.02
$ws.Cells.Item(15, 3).Value2 = 1.05642505648213
$ws.Cells.Item(15, 4).Value2 = 1.057244970534024
$ws.Cells.Item(15, 5).Value2 = 1.061829307361001
$ws.Cells.Item(15, 6).Value2 = 1.070344029494678
$ws.Cells.Item(15, 9).Value2 = 1.040345352515204
$ws.Cells.Item(15, 10).Value2 = 1.062814008863874
$ws.Cells.Item(15, 11).Value2 = 1.060719645869813
$ws.Cells.Item(15, 12).Value2 = 1.065287886225752
$ws.Cells.Item(15, 13).Value2 = 1.073773114196776
$ws.Cells.Item(15, 14).Value2 = 1.064323326662932

# Row 16
$ws.Cells.Item(16, 2).Value2 = 1.02
$ws.Cells.Item(16, 3).Value2 = 1.057373370892355
$ws.Cells.Item(16, 4).Value2 = 1.057962378983943
$ws.Cells.Item(16, 5).Value2 = 1.062658578408304
$ws.Cells.Item(16, 6).Value2 = 1.071193211967866
$ws.Cells.Item(16, 9).Value2 = 1.040542828599373
$ws.Cells.Item(16, 10).Value2 = 1.063546498775422
$ws.Cells.Item(16, 11).Value2 = 1.061323406139063
$ws.Cells.Item(16, 12).Value2 = 1.066003701712233
$ws.Cells.Item(16, 13).Value2 = 1.074509824238914
$ws.Cells.Item(16, 14).Value2 = 1.065056856794172

# Row 17
$ws.Cells.Item(17, 2).Value2 = 1.02
$ws.Cells.Item(17, 3).Value2 = 1.057967825224957
$ws.Cells.Item(17, 4).Value2 = 1.058411988198012
$ws.Cells.Item(17, 5).Value2 = 1.063178484012918
$ws.Cells.Item(17, 6).Value2 = 1.071725610206827
$ws.Cells.Item(17, 9).Value2 = 1.040666082473911
$ws.Cells.Item(17, 10).Value2 = 1.064005452676756
$ws.Cells.Item(17, 11).Value2 = 1.061701538899824
$ws.Cells.Item(17, 12).Value2 = 1.066452264857682
$ws.Cells.Item(17, 13).Value2 = 1.074971497961776
$ws.Cells.Item(17, 14).Value2 = 1.065516462462688

# Row 18
$ws.Cells.Item(18, 2).Value2 = 1.02
$ws.Cells.Item(18, 3).Value2 = 1.058314419390383
$ws.Cells.Item(18, 4).Value2 = 1.058674094278099
$ws.Cells.Item(18, 5).Value2 = 1.063481639208103
$ws.Cells.Item(18, 6).Value2 = 1.072036053128564
$ws.Cells.Item(18, 9).Value2 = 1.040737752339195
$ws.Cells.Item(18, 10).Value2 = 1.064272968001569
$ws.Cells.Item(18, 11).Value2 = 1.061921886149448
$ws.Cells.Item(18, 12).Value2 = 1.06671374406588
$ws.Cells.Item(18, 13).Value2 = 1.07524062577867
$ws.Cells.Item(18, 14).Value2 = 1.065784357689947

# Row 19
$ws.Cells.Item(19, 2).Value2 = 1.02
$ws.Cells.Item(19, 3).Value2 = 1.058432575366763
$ws.Cells.Item(19, 4).Value2 = 1.058763441506365
$ws.Cells.Item(19, 5).Value2 = 1.063584991026755
$ws.Cells.Item(19, 6).Value2 = 1.072141890040848
$ws.Cells.Item(19, 9).Value2 = 1.040762152302536
$ws.Cells.Item(19, 10).Value2 = 1.064364152654296
$ws.Cells.Item(19, 11).Value2 = 1.061996983183155
$ws.Cells.Item(19, 12).Value2 = 1.066802874726898
$ws.Cells.Item(19, 13).Value2 = 1.075332364660077
$ws.Cells.Item(19, 14).Value2 = 1.06587567183534

# Row 20
$ws.Cells.Item(20, 2).Value2 = 1.02
$ws.Cells.Item(20, 3).Value2 = 1.05790406054124
$ws.Cells.Item(20, 4).Value2 = 1.058363764248891
$ws.Cells.Item(20, 5).Value2 = 1.063122713131995
$ws.Cells.Item(20, 6).Value2 = 1.07166849888273
$ws.Cells.Item(20, 9).Value2 = 1.040652881478356
$ws.Cells.Item(20, 10).Value2 = 1.063956230414713
$ws.Cells.Item(20, 11).Value2 = 1.06166099071108
$ws.Cells.Item(20, 12).Value2 = 1.066404154874713
$ws.Cells.Item(20, 13).Value2 = 1.074921981189912
$ws.Cells.Item(20, 14).Value2 = 1.065467170299387

# Row 21
$ws.Cells.Item(21, 2).Value2 = 1.02
$ws.Cells.Item(21, 3).Value2 = 1.056184132909335
$ws.Cells.Item(21, 4).Value2 = 1.057062678849176
$ws.Cells.Item(21, 5).Value2 = 1.061618650400233
$ws.Cells.Item(21, 6).Value2 = 1.070128317404966
$ws.Cells.Item(21, 9).Value2 = 1.040295017911566
$ws.Cells.Item(21, 10).Value2 = 1.062627851805763
$ws.Cells.Item(21, 11).Value2 = 1.060566154173584
$ws.Cells.Item(21, 12).Value2 = 1.065105984451738
$ws.Cells.Item(21, 13).Value2 = 1.073585908044506
$ws.Cells.Item(21, 14).Value2 = 1.064136905240446

# Row 22
$ws.Cells.Item(22, 2).Value2 = 1.02
$ws.Cells.Item(22, 3).Value2 = 1.055101421689689
$ws.Cells.Item(22, 4).Value2 = 1.056243309712649
$ws.Cells.Item(22, 5).Value2 = 1.060672072854941
$ws.Cells.Item(22, 6).Value2 = 1.069159039373735
$ws.Cells.Item(22, 9).Value2 = 1.040068002573877
$ws.Cells.Item(22, 10).Value2 = 1.061790943864421
$ws.Cells.Item(22, 11).Value2 = 1.059875852769685
$ws.Cells.Item(22, 12).Value2 = 1.064288293945254
$ws.Cells.Item(22, 13).Value2 = 1.072744399517408
$ws.Cells.Item(22, 14).Value2 = 1.063298808793832

# Row 23
$ws.Cells.Item(23, 2).Value2 = 1.02
$ws.Cells.Item(23, 3).Value2 = 1.055675515575282
$ws.Cells.Item(23, 4).Value2 = 1.056677800029796
$ws.Cells.Item(23, 5).Value2 = 1.061173960076737
$ws.Cells.Item(23, 6).Value2 = 1.069672959721348
$ws.Cells.Item(23, 9).Value2 = 1.040188539027556
$ws.Cells.Item(23, 10).Value2 = 1.062234767997465
$ws.Cells.Item(23, 11).Value2 = 1.060241979483249
$ws.Cells.Item(23, 12).Value2 = 1.064721909169626
$ws.Cells.Item(23, 13).Value2 = 1.073190639898294
$ws.Cells.Item(23, 14).Value2 = 1.063743263208052

# Row 24
$ws.Cells.Item(24, 2).Value2 = 1.02
$ws.Cells.Item(24, 3).Value2 = 1.057932873493054
$ws.Cells.Item(24, 4).Value2 = 1.058385555019054
$ws.Cells.Item(24, 5).Value2 = 1.063147913892364
$ws.Cells.Item(24, 6).Value2 = 1.071694305328646
$ws.Cells.Item(24, 9).Value2 = 1.040658847126211
$ws.Cells.Item(24, 10).Value2 = 1.063978472407254
$ws.Cells.Item(24, 11).Value2 = 1.061679313345041
$ws.Cells.Item(24, 12).Value2 = 1.066425894200099
$ws.Cells.Item(24, 13).Value2 = 1.074944356178244
$ws.Cells.Item(24, 14).Value2 = 1.065489443878109

# Row 25
$ws.Cells.Item(25, 2).Value2 = 1.02
$ws.Cells.Item(25, 3).Value2 = 1.060546908345534
$ws.Cells.Item(25, 4).Value2 = 1.060361702393356
$ws.Cells.Item(25, 5).Value2 = 1.065434794558297
$ws.Cells.Item(25, 6).Value2 = 1.074036220166927
$ws.Cells.Item(25, 9).Value2 = 1.041195909520639
$ws.Cells.Item(25, 10).Value2 = 1.065994717912388
$ws.Cells.Item(25, 11).Value2 = 1.063338990484928
$ws.Cells.Item(25, 12).Value2 = 1.068397012757731
$ws.Cells.Item(25, 13).Value2 = 1.076973241177083
$ws.Cells.Item(25, 14).Value2 = 1.067508552683126
